$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "clean" numeric-looking string must be forced to
# text (NumberFormat "@") before assignment, otherwise Excel auto-converts them
# to numbers -- the source data keeps these as plain text cells.

$ws.Range("D2").Value = "68.557.97"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "2.655.65"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.90"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.39"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").Value = "2.656.77"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.146"
$ws.Range("E10").Value = "  +8.75%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.38"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000194"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "3.130.56"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "68.373.75"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "2.669.29"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.57"
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "367.66"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.57"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.47"
$ws.Range("E22").Value = "  +5.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.99"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000109"
$ws.Range("E28").Value = "  +3.75%  "
$ws.Range("D29").Value = "2.773.95"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "584.49"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.19"
$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.133"
$ws.Range("E35").Value = "  +4.57%  "
$ws.Range("E36").Value = "  +5.67%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.76"
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.57"
$ws.Range("E39").Value = "  +2.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.91"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.376"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.46"
$ws.Range("E42").Value = "  +3.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.72"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₆0332"
$ws.Range("E44").Value = "  +11.95%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.73"
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.53"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "158.36"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.81"
$ws.Range("E49").Value = "  +3.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.73"
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.14"
$ws.Range("E51").Value = "  +3.45%  "
